$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

function Set-PlainCell($addr, $val) {
    $ws.Range($addr).Value = $val
}

Set-TextCell "D2" "66.285.82"
Set-PlainCell "E2" "  -1.81%  "
Set-TextCell "D3" "3.515.45"
Set-PlainCell "E3" "  -0.84%  "
Set-TextCell "D4" "0.999"
Set-PlainCell "E4" "  -0.18%  "
Set-TextCell "D5" "584.27"
Set-PlainCell "E5" "  +5.72%  "
Set-TextCell "D6" "179.49"
Set-PlainCell "E6" "  -5.49%  "
Set-TextCell "D7" "0.638"
Set-PlainCell "E7" "  +4.73%  "
Set-PlainCell "E8" "  -0.02%  "
Set-TextCell "D9" "0.642"
Set-PlainCell "E9" "  +1.14%  "
Set-TextCell "D10" "0.164"
Set-PlainCell "E10" "  +5.94%  "
Set-TextCell "D11" "56.25"
Set-PlainCell "E11" "  +2.41%  "
Set-TextCell "D12" "0.0000281"
Set-PlainCell "E12" "  +3.52%  "
Set-TextCell "D13" "9.34"
Set-PlainCell "E13" "  -0.69%  "
Set-TextCell "D14" "4.078.17"
Set-PlainCell "E14" "  -0.60%  "
Set-TextCell "D15" "3.518.65"
Set-PlainCell "E15" "  -0.53%  "
Set-PlainCell "E16" "  +0.06%  "
Set-TextCell "D17" "18.45"
Set-PlainCell "E17" "  +1.04%  "
Set-TextCell "D18" "66.282.37"
Set-PlainCell "E18" "  -1.86%  "
Set-TextCell "D19" "12.07"
Set-PlainCell "E20" "  +2.15%  "
Set-TextCell "D21" "415.72"
Set-PlainCell "E21" "  -3.83%  "
Set-PlainCell "E22" "  +9.69%  "
Set-TextCell "D23" "4.43"
Set-PlainCell "E23" "  +6.80%  "
Set-TextCell "D24" "85.37"
Set-PlainCell "E24" "  -0.17%  "
Set-TextCell "D25" "13.51"
Set-PlainCell "E25" "  +11.56%  "
Set-TextCell "D26" "11.13"
Set-PlainCell "E26" "  -0.02%  "
Set-PlainCell "E27" "  -1.55%  "
Set-TextCell "D28" "6.05"
Set-PlainCell "E28" "  -1.69%  "
Set-PlainCell "E29" "  +2.11%  "
Set-TextCell "D30" "30.51"
Set-PlainCell "E30" "  +0.37%  "
Set-TextCell "D31" "6.74"
Set-PlainCell "E31" "  +0.48%  "
Set-PlainCell "B32" "Cosmos"
Set-PlainCell "C32" "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextCell "D32" "11.82"
Set-PlainCell "E32" "  +0.46%  "
Set-PlainCell "B33" "Bittensor"
Set-PlainCell "C33" "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextCell "D33" "598.04"
Set-PlainCell "E33" "  -7.25%  "
Set-TextCell "D34" "0.112"
Set-PlainCell "E34" "  +0.12%  "
Set-TextCell "D35" "60.90"
Set-PlainCell "E35" "  +1.69%  "
Set-PlainCell "E36" "  +5.46%  "
Set-TextCell "D37" "0.0₃0804"
Set-PlainCell "E37" "  -3.06%  "
Set-TextCell "D38" "1.00"
Set-PlainCell "E38" "  +0.11%  "
Set-TextCell "D39" "3.66"
Set-PlainCell "E39" "  +9.09%  "
Set-TextCell "D40" "36.98"
Set-PlainCell "E40" "  -4.31%  "
Set-PlainCell "E41" "  -1.27%  "
Set-TextCell "D42" "3.240.24"
Set-PlainCell "E42" "  +6.48%  "
Set-TextCell "D43" "0.999"
Set-PlainCell "E43" "  -0.11%  "
Set-PlainCell "E44" "  +2.94%  "
Set-PlainCell "B45" "ApeXProtocol"
Set-PlainCell "C45" "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
Set-TextCell "D45" "3.35"
Set-PlainCell "E45" "  +0.95%  "
Set-PlainCell "B46" "Fetch.AI"
Set-PlainCell "C46" "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextCell "D46" "2.57"
Set-PlainCell "E46" "  -3.43%  "
Set-PlainCell "E47" "  +1.07%  "
Set-PlainCell "E48" "  +1.93%  "
Set-PlainCell "E49" "  -5.90%  "
Set-PlainCell "E50" "  -0.80%  "
Set-TextCell "D51" "140.16"
Set-PlainCell "E51" "  -0.66%  "

Write-Host "Applied all changes"
